# Updates cryptos list price (D) and volume-change (E) columns per the
# refreshed coinranking.com snapshot (GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.505.60"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.907.17"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.338"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0705"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "2.187.89"
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.68%  "
$ws.Range("D14").Value = "1.916.87"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.691"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").Value = "35.548.91"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "72.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "0.0₃0827"
$ws.Range("E19").Value = "  +4.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.964"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +24.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0570"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("E35").Value = "  +6.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "91.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.03%  "
$ws.Range("D42").Value = "1.352.50"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  +13.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "49.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +42.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").Value = "2.096.52"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0691"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.13%  "
